$d = $word.ActiveDocument

# --- Paragraph 1 formatting: add a paragraph border (space-only, no line)
# and widen the left indent from 120 twips (6pt) to 225 twips (11.25pt),
# matching the pBdr/indent already used further down in the document.
$p1 = $d.Paragraphs(1)
$p1.Range.Borders.DistanceFromTop = 5
$p1.Range.Borders.DistanceFromLeft = 5
$p1.Range.Borders.DistanceFromBottom = 5
$p1.Range.Borders.DistanceFromRight = 5
$p1.LeftIndent = 11.25

# --- Update the placeholder id text and drop the now-unwanted trailing
# space run that used to follow it (same rPr, so the engine folds them
# into a single run automatically once the text matches).
$d.Content.Find.Execute("**ID__AFFARS_pgi_5305_topic_8__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_SMC_PGI_5305_303__ID**", 2)
$d.Content.Find.Execute("**ID__AFFARS_SMC_PGI_5305_303__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_SMC_PGI_5305_303__ID**", 2)

Write-Host "Done"
